{"js": "// Apply the \"Forming the Community\" MC-script restructuring + date fix.\n\nfunction findByText(items, text) {\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text === text) return items[i];\n  }\n  throw new Error(\"Paragraph not found: \" + text);\n}\n\n// --- 1) Fix the date: \"May 24, 2018\" -> \"May 28, 2018\" ---\nconst dateSearch = context.document.body.search(\"24,\", { matchCase: true });\ndateSearch.load(\"items\");\nawait context.sync();\nif (dateSearch.items.length > 0) {\n  dateSearch.items[0].insertText(\"28,\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 2) Restructure the paragraphs under \"Forming the Community\" ---\nlet paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet items = paragraphs.items;\n\nconst goalPara = findByText(\n  items,\n  \"Goal: Participants develop familiiarty with each other, build a strong sense of unity, and common purpose, and a willingness to be active participants in the various learning experiences of the LDZ.\"\n);\nconst welcomeHeadingPara = findByText(items, \"OFFICIAL WELCOME\");\n\n// 2a) \"Goal: ...\" paragraph becomes a Heading 4 with a bookmark wrapping its (zero-length) start.\ngoalPara.style = \"Heading 4\";\n\n// 2b) Insert a new \"Info: MC is a staff member, chosen by ED's\" Heading 4 paragraph after \"OFFICIAL WELCOME\".\nconst infoPara = welcomeHeadingPara.insertParagraph(\n  \"Info: MC is a staff member, chosen by ED\\u2019s\",\n  Word.InsertLocation.after\n);\ninfoPara.style = \"Heading 4\";\n\n// Sync so the newly inserted paragraph/style changes settle before we touch ranges/bookmarks\n// (bookmarking a just-inserted-and-not-yet-synced paragraph's Start range can anchor to the\n// wrong paragraph, so re-fetch fresh references after the sync below).\nawait context.sync();\n\nparagraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\nitems = paragraphs.items;\n\nconst goalPara2 = findByText(\n  items,\n  \"Goal: Participants develop familiiarty with each other, build a strong sense of unity, and common purpose, and a willingness to be active participants in the various learning experiences of the LDZ.\"\n);\ngoalPara2\n  .getRange(\"Start\")\n  .insertBookmark(\n    \"goal-participants-develop-familiiarty-with-each-other-build-a-strong-sense-of-unity-and-common-purpose-and-a-willingness-to-be-active-participants-in-the-various-learning-experiences-of-the-ldz.\"\n  );\n\nconst infoPara2 = findByText(items, \"Info: MC is a staff member, chosen by ED\\u2019s\");\ninfoPara2.getRange(\"Start\").insertBookmark(\"info-mc-is-a-staff-member-chosen-by-eds\");\n\nconst deliveredPara = findByText(items, \"(DELIVERED BY A STAFF MEMBER - Called the MC)\");\nconst emphasisPara = findByText(\n  items,\n  \"[With Emphasis] Welcome DELEGATES to the 2018\\u2026NATIONAL\\u2026Lorenzo de Zavala\\u2026Youth\\u2026Legislative\\u2026Session\\u2026\"\n);\nconst rollCallPara = findByText(items, \"We will now begin the official roll call. Who is here from\\u2026\");\nconst thankYouPara = findByText(items, \"Thank you all for being with us, we will continue with our program\");\nconst welcomeEdPara = findByText(items, \"I would now like to welcome to the stage, Education Director:\");\n\n// 2c) \"(DELIVERED BY A STAFF MEMBER - Called the MC)\" -> new MC welcome line (style unchanged: First Paragraph).\ndeliveredPara.insertText(\n  \"MC: Welcome DELEGATES to the 2018 - NATIONAL - Lorenzo de Zavala - Youth - Legislative - Session!!!!!!!!!!\",\n  Word.InsertLocation.replace\n);\n\n// 2d) \"[With Emphasis] Welcome DELEGATES...\" -> \"MC: We will now begin the official roll call. Who is here from...\"\nemphasisPara.insertText(\n  \"MC: We will now begin the official roll call. Who is here from\\u2026\",\n  Word.InsertLocation.replace\n);\n\n// 2e) \"We will now begin the official roll call...\" -> merged \"Thank you\"/\"welcome ED\" MC line.\nrollCallPara.insertText(\n  \"MC: Thank you all for being with us, we will now continue with our program. I would now like to welcome to the stage, Education Director:\",\n  Word.InsertLocation.replace\n);\n\n// 2f) Delete the now-redundant trailing paragraphs.\nthankYouPara.delete();\nwelcomeEdPara.delete();\n\nawait context.sync();\n", "ps1": "# Apply the \"Forming the Community\" MC-script restructuring + date fix.\n\n$d = $word.ActiveDocument\n\n# Special (non-ASCII) characters used by the source text, built from code points so the\n# interpreter never has to round-trip literal UTF-8 through the script text itself.\n$ellipsis = [char]0x2026\n$rsquo = [char]0x2019\n\n# --- 1) Fix the date: \"May 24, 2018\" -> \"May 28, 2018\" ---\n$rng = $d.Content\n[void]$rng.Find.Execute(\"24,\", $false, $false, $false, $false, $false, $true, 1, $false, \"28,\", 2)\n\n# --- Helper: find a paragraph whose text equals $text (trimming the trailing paragraph mark). ---\nfunction Get-ParagraphByText($doc, $text) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        $p = $doc.Paragraphs($i)\n        $t = $p.Range.Text\n        $t = $t.TrimEnd([char]13, [char]7)\n        if ($t -eq $text) {\n            return $p\n        }\n    }\n    return $null\n}\n\n# --- 2) Restructure the paragraphs under \"Forming the Community\" ---\n\n# 2a) \"Goal: ...\" paragraph becomes a Heading 4 with a bookmark wrapping its (zero-length) start.\n$goalPara = Get-ParagraphByText $d \"Goal: Participants develop familiiarty with each other, build a strong sense of unity, and common purpose, and a willingness to be active participants in the various learning experiences of the LDZ.\"\n$goalPara.Style = \"Heading 4\"\n$s = $goalPara.Range.Start\n$collapsed = $d.Range($s, $s)\n$d.Bookmarks.Add(\"goal-participants-develop-familiiarty-with-each-other-build-a-strong-sense-of-unity-and-common-purpose-and-a-willingness-to-be-active-participants-in-the-various-learning-experiences-of-the-ldz.\", $collapsed)\n\n# 2b) Insert a new \"Info: MC is a staff member, chosen by ED's\" Heading 4 paragraph after \"OFFICIAL WELCOME\".\n$welcomeHeadingPara = Get-ParagraphByText $d \"OFFICIAL WELCOME\"\n$welcomeHeadingPara.Range.InsertParagraphAfter()\n$infoPara = $welcomeHeadingPara.Next()\n$infoPara.Range.Text = \"Info: MC is a staff member, chosen by ED\" + $rsquo + \"s\"\n$infoPara.Style = \"Heading 4\"\n$s2 = $infoPara.Range.Start\n$collapsed2 = $d.Range($s2, $s2)\n$d.Bookmarks.Add(\"info-mc-is-a-staff-member-chosen-by-eds\", $collapsed2)\n\n# 2c) \"(DELIVERED BY A STAFF MEMBER - Called the MC)\" -> new MC welcome line (style unchanged: First Paragraph).\n$deliveredPara = Get-ParagraphByText $d \"(DELIVERED BY A STAFF MEMBER - Called the MC)\"\n$deliveredPara.Range.Text = \"MC: Welcome DELEGATES to the 2018 - NATIONAL - Lorenzo de Zavala - Youth - Legislative - Session!!!!!!!!!!\"\n\n# 2d) \"[With Emphasis] Welcome DELEGATES...\" -> \"MC: We will now begin the official roll call. Who is here from...\"\n$emphasisText = \"[With Emphasis] Welcome DELEGATES to the 2018\" + $ellipsis + \"NATIONAL\" + $ellipsis + \"Lorenzo de Zavala\" + $ellipsis + \"Youth\" + $ellipsis + \"Legislative\" + $ellipsis + \"Session\" + $ellipsis\n$emphasisPara = Get-ParagraphByText $d $emphasisText\n$emphasisPara.Range.Text = \"MC: We will now begin the official roll call. Who is here from\" + $ellipsis\n\n# 2e) \"We will now begin the official roll call...\" -> merged \"Thank you\"/\"welcome ED\" MC line.\n$rollCallText = \"We will now begin the official roll call. Who is here from\" + $ellipsis\n$rollCallPara = Get-ParagraphByText $d $rollCallText\n$rollCallPara.Range.Text = \"MC: Thank you all for being with us, we will now continue with our program. I would now like to welcome to the stage, Education Director:\"\n\n# 2f) Delete the now-redundant trailing paragraphs.\n$thankYouPara = Get-ParagraphByText $d \"Thank you all for being with us, we will continue with our program\"\n[void]$thankYouPara.Range.Delete()\n\n$welcomeEdPara = Get-ParagraphByText $d \"I would now like to welcome to the stage, Education Director:\"\n[void]$welcomeEdPara.Range.Delete()\n"}
